$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.329.02"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "1.842.13"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9984"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6300"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07520"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2938"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.40"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07687"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("D12").Value = "1.836.19"
$ws.Range("E12").Value = "  -7.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.979"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6775"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001051"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.92"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "2.085.39"
$ws.Range("E17").Value = "  -7.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.116"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "29.367.71"
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "227.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.62%  "
$ws.Range("E21").Value = "  -0.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9994"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.411"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.61%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E25").Value = "  +1.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1386"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.346"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("E29").Value = "  -1.19%  "
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05616"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.099"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.015"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.827"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.154"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7082"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.34%  "
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").Value = "1.239.97"
$ws.Range("E38").Value = "  -0.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01807"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.764"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.237"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9010"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9990"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000121"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.074"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.3989"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("E51").Value = "  -0.39%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.898"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.85%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.667"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.43%  "
